$d = $word.ActiveDocument

# --- paragraph 36 (p36) ---
$r = $d.Paragraphs(36).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">The project has been delegated to team members in the form of design and function coding. As such, all the formatting and graphic functions have been reviewed by each of the team members during meetings and are walked through and altered to make sure the program is written in an effective manner and to ensure that a team member does not do any work that another team member </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>can not</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> integrate.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 34 (p34) ---
$r = $d.Paragraphs(34).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="288" w:firstLine="432"/><w:rPr><w:rFonts w:ascii="Merriweather Sans" w:eastAsia="Merriweather Sans" w:hAnsi="Merriweather Sans" w:cs="Merriweather Sans"/><w:i/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Execution-based Non-Functional Testing is done by all the team members trying to log in into the app and see if they can login when they do not have an internet connection as the program require strong internet connection to run properly. </w:t></w:r><w:r><w:t xml:space="preserve">Also checked if every member can run the </w:t></w:r><w:r><w:t xml:space="preserve">map with specified location </w:t></w:r><w:r><w:t xml:space="preserve">and it can </w:t></w:r><w:r><w:t xml:space="preserve">run on all devices without breaking the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mapview</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Made sure all members could push and pull from the firebase regardless of system. Also made sure that the app was displaying properly on </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>various different</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> platforms and operating system versions, tested android </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oreo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and pie, current windows </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>linux</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> distributions.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 32 (p32) ---
$r = $d.Paragraphs(32).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Functional testing of the project was done by each group member. After each push to the repository group member compiled and ran the application on their devices. Testers tested opening the app to see if the app open or crashed. Functionalities such as moving through tabs was also tested. Functional testing </w:t></w:r><w:r><w:t xml:space="preserve">for login page is done by each member, to see if we can </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>signup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the system and login and see if the details were saved in the database. Every member ran the map view function added to the app and test if it is accurate for the geotags provided.</w:t></w:r><w:r><w:t xml:space="preserve"> Also working on connecting the actual google maps for android users. All these functionalities are tested by the members.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 29 (p29) ---
$r = $d.Paragraphs(29).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>GeoPy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Used to permit applications Geo tagging system</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 28 (p28) ---
$r = $d.Paragraphs(28).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Firebase and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pyrebase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Used for storing and retrieving login information.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 27 (p27) ---
$r = $d.Paragraphs(27).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Kivy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Used throughout the entire project.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 25 (p25) ---
$r = $d.Paragraphs(25).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Python is used throughout the entire project as well as </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Kivy’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> format for the layouts and string assignments. Worth noting that any,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r><w:t>kv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> text file functions as a sort of Strings.txt file in Android Studio projects.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

# --- paragraph 16 (p16) ---
$r = $d.Paragraphs(16).Range
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Brisan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> B</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)
